$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $find"
    }
}

# Question 1: "camera controls" wording change
Replace-Text "how satisfying was the controls? If unsatisfying" "how satisfying are the camera controls? If unsatisfying"

# Question 3: replace the ratings-scale qualifier clause and trailing question
Replace-Text "Question 3: On a scale of 1 -5(from 1 being least satisfying, to 5 being most satisfying), how smooth did you find the movement of the player? If a 1, why was it unsatisfactory?" "Question 3: On a scale of 1 -5(from 1 being least satisfying) what was your impression of the UI and how would you improve it?"

# Question 5: new wording
Replace-Text "Question 5: Which mechanic or feature of the game did you enjoy the most and why?" "Question 5: How clear are the instructions and what else would you include if it was unclear?"

# Question 7: new wording
Replace-Text "Question 7: What was your impression on the map generation and do you find it effective on a scale of 1 – 5, with 5 being very effective?" "Question 7: What is your impression of the basic combat animation? Would you suggest any improvements?"

# Question 8: new wording
Replace-Text "Question 8: What features did you find unnecessary, distracting or was completely ignored when playing the game?" "Question 8: Based on the Current UI & Controls, is there any changes or issues with how you interact with the game?"

# Question 10: new wording (prototype -> game)
Replace-Text "what was your impression of the prototype and why?" "what was your impression of the game and why?"

# --- Ratings line after Question 7: clear the "1 2 3 4 5" line and un-center it ---
$q7ratings = $d.Paragraphs.Item(14)
$delRange = $d.Range($q7ratings.Range.Start, $q7ratings.Range.End - 1)
$delRange.Text = ""
$q7ratings.Alignment = 0

# --- Ratings line after Question 10: move the _GoBack bookmark from the end of the
#     line to right after the "3" (before the following space) ---
$q10ratings = $d.Paragraphs.Item(18)
$bmPos = $q10ratings.Range.Start + 6
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Done. Final paragraph text:"
for ($i = 1; $i -le 18; $i++) {
    Write-Host $i ":" $d.Paragraphs.Item($i).Range.Text
}
